$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Rows 3-19 hold the training records. Add a new day of progress:
#  - col H ("PERIOD TO EXPIRE") drops by one day
#  - col I ("LAST UPDATE") moves from 03-Nov-2025 to 04-Nov-2025
# Force column I to stay plain text so "04-Nov-2025" is not auto-converted
# into a date serial value by Excel's type inference.
$lastUpdateRange = $ws.Range("I3:I19")
$lastUpdateRange.NumberFormat = "@"

for ($row = 3; $row -le 19; $row++) {
    $periodCell = $ws.Cells.Item($row, 8)
    $periodCell.Value2 = $periodCell.Value2 - 1

    $ws.Cells.Item($row, 9).Value2 = "04-Nov-2025"
}
